# Auto-generated script applying scheduled-runner cell updates to Sheets/Ravana_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 998.8333
$ws.Range("I2").Value = 998.8333
$ws.Range("K2").Value = 998.8333
$ws.Range("M2").Value = -885.8333
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H34").Value = 7044
$ws.Range("I34").Value = 7044
$ws.Range("K34").Value = 7044
$ws.Range("M34").Value = -6841
$ws.Range("H36").Value = 7044
$ws.Range("I36").Value = 7044
$ws.Range("K36").Value = 7044
$ws.Range("M36").Value = -6329
$ws.Range("H64").Value = 3257.2856
$ws.Range("J64").Value = 3280
$ws.Range("L64").Value = 3280
$ws.Range("N64").Value = -3776
$ws.Range("H67").Value = 3257.2856
$ws.Range("J67").Value = 3280
$ws.Range("L67").Value = 3280
$ws.Range("N67").Value = -4996
$ws.Range("H69").Value = 3166.3333
$ws.Range("I69").Value = 3166.3333
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 9498.999899999999
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -8624.999899999999
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 3166.3333
$ws.Range("I72").Value = 3166.3333
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 28496.9997
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -24128.9997
$ws.Range("N72").ClearContents()
$ws.Range("H132").Value = 1288.9333
$ws.Range("I132").Value = 1288.9333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3866.7999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1336.7999
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 75000
$ws.Range("J133").Value = 75000
$ws.Range("L133").Value = 75000
$ws.Range("N133").Value = -85120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H23").Value = 35997
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 35997
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 35997
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -36515
$ws.Range("H61").Value = 2124.1
$ws.Range("I61").Value = 1548.6
$ws.Range("J61").Value = 2699.6
$ws.Range("K61").Value = 1548.6
$ws.Range("L61").Value = 2699.6
$ws.Range("M61").Value = -1336.6
$ws.Range("N61").Value = -3123.6
$ws.Range("H132").Value = 2525.05
$ws.Range("I132").Value = 1583.8334
$ws.Range("J132").Value = 3936.875
$ws.Range("K132").Value = 4751.5002
$ws.Range("L132").Value = 11810.625
$ws.Range("M132").Value = -2221.5002
$ws.Range("N132").Value = -16870.625
$ws.Range("H136").Value = 2124.1
$ws.Range("I136").Value = 1548.6
$ws.Range("J136").Value = 2699.6
$ws.Range("K136").Value = 4645.799999999999
$ws.Range("L136").Value = 8098.799999999999
$ws.Range("M136").Value = -2095.799999999999
$ws.Range("N136").Value = -13198.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 676.8
$ws.Range("I22").Value = 656
$ws.Range("J22").Value = 760
$ws.Range("K22").Value = 656
$ws.Range("L22").Value = 760
$ws.Range("M22").Value = -483
$ws.Range("N22").Value = -1106
$ws.Range("H64").Value = 1316
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1316
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 1316
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -1766
$ws.Range("H67").Value = 1316
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1316
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 1316
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -2876
$ws.Range("H134").Value = 2757.3333
$ws.Range("I134").Value = 2714.5454
$ws.Range("K134").Value = 8143.6362
$ws.Range("M134").Value = -5608.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2744.5
$ws.Range("I58").Value = 2449
$ws.Range("J58").Value = 2980.9
$ws.Range("K58").Value = 2449
$ws.Range("L58").Value = 2980.9
$ws.Range("M58").Value = -2246
$ws.Range("N58").Value = -3386.9
$ws.Range("H107").Value = 2546.7
$ws.Range("I107").Value = 935
$ws.Range("K107").Value = 935
$ws.Range("M107").Value = 985
$ws.Range("H132").Value = 4999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 4000
$ws.Range("I134").Value = 4000
$ws.Range("K134").Value = 12000
$ws.Range("M134").Value = -9465
$ws.Range("H136").Value = 2744.5
$ws.Range("I136").Value = 2449
$ws.Range("J136").Value = 2980.9
$ws.Range("K136").Value = 7347
$ws.Range("L136").Value = 8942.700000000001
$ws.Range("M136").Value = -4797
$ws.Range("N136").Value = -14042.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2561.375
$ws.Range("J113").Value = 2561.375
$ws.Range("L113").Value = 7684.125
$ws.Range("N113").Value = -12024.125
$ws.Range("H140").Value = 1401.6666
$ws.Range("I140").Value = 1401.6666
$ws.Range("K140").Value = 4204.9998
$ws.Range("M140").Value = 975.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 111.666664
$ws.Range("I2").Value = 105
$ws.Range("J2").Value = 125
$ws.Range("K2").Value = 105
$ws.Range("L2").Value = 125
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = -351
$ws.Range("H80").Value = 5162.25
$ws.Range("I80").Value = 4470.2856
$ws.Range("J80").Value = 10006
$ws.Range("K80").Value = 4470.2856
$ws.Range("L80").Value = 10006
$ws.Range("M80").Value = -3472.2856
$ws.Range("N80").Value = -12002
$ws.Range("H83").Value = 5162.25
$ws.Range("I83").Value = 4470.2856
$ws.Range("J83").Value = 10006
$ws.Range("K83").Value = 22351.428
$ws.Range("L83").Value = 50030
$ws.Range("M83").Value = -17359.428
$ws.Range("N83").Value = -60014
$ws.Range("H97").Value = 405.14285
$ws.Range("I97").Value = 444.5909
$ws.Range("J97").Value = 260.5
$ws.Range("K97").Value = 444.5909
$ws.Range("L97").Value = 260.5
$ws.Range("M97").Value = 51.40910000000002
$ws.Range("N97").Value = -1252.5
$ws.Range("H102").Value = 2083.3333
$ws.Range("I102").Value = 1625
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1625
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -3
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 3426.2727
$ws.Range("I132").Value = 3138.2
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 9414.599999999999
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -6884.599999999999
$ws.Range("N132").Value = -16058.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5120.364
$ws.Range("I22").Value = 3724.8
$ws.Range("J22").Value = 6283.3335
$ws.Range("K22").Value = 3724.8
$ws.Range("L22").Value = 6283.3335
$ws.Range("M22").Value = -3429.8
$ws.Range("N22").Value = -6873.3335
$ws.Range("H27").Value = 5120.364
$ws.Range("I27").Value = 3724.8
$ws.Range("J27").Value = 6283.3335
$ws.Range("K27").Value = 3724.8
$ws.Range("L27").Value = 6283.3335
$ws.Range("M27").Value = -3617.8
$ws.Range("N27").Value = -6497.3335
$ws.Range("H43").Value = 459992.66
$ws.Range("J43").Value = 459992.66
$ws.Range("L43").Value = 459992.66
$ws.Range("N43").Value = -460378.66
$ws.Range("H132").Value = 5311.75
$ws.Range("I132").Value = 5082.6665
$ws.Range("K132").Value = 15247.9995
$ws.Range("M132").Value = -12717.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 593.7778
$ws.Range("I107").Value = 540.6667
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 1622.0001
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = 297.9999
$ws.Range("N107").Value = -5940
$ws.Range("H122").Value = 4528.143
$ws.Range("I122").Value = 4449.5
$ws.Range("K122").Value = 13348.5
$ws.Range("M122").Value = -10898.5
